# Mainboard: update BOM, fabrication print
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mainboard_bom_")
$ws.Activate()

# --- Update BOM LCSC part numbers ---
# Row 12 (USB connector): LCSC code corrected
$ws.Range("J12").Value = "C132563"

# Row 27 (Female header, 14p, 0.1"): add LCSC code
$ws.Range("J27").Value = "C52711"

# Row 29 (Female header, 2*4p, 0.1"): add LCSC code
$ws.Range("J29").Value = "C92271"

# Row 30 (Female header, 2*6p, 0.1"): add LCSC code
$ws.Range("J30").Value = "C92269"

# --- Update view/window state to reflect author's scroll position & selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Rows(27).Select()
